$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# Column D header changes from "Adresse" to "Straße"
$ws.Range("D1").Value = "Straße"

# New columns H (E-Nummer) and I (Befundtyp)
$ws.Range("H1").Value = "E-Nummer"
$ws.Range("I1").Value = "Befundtyp"

# --- Data rows ---
$ws.Range("H2").Value = "A/1996/200592"
$ws.Range("I2").Value = "Hauptbefund"

$ws.Range("H3").Value = "A/1996/200591"
$ws.Range("I3").Value = "Hauptbefund"

$ws.Range("H4").Value = "A/1996/200391"
$ws.Range("I4").Value = "Nebenbefund"

$ws.Range("H5").Value = "A/1998/200591"
$ws.Range("I5").Value = "Hauptbefund"

$ws.Range("H6").Value = "001/00146"
$ws.Range("I6").Value = "Hauptbefund"

# H5, H6 and H7 carry the date-format style (s="1") inherited from column E
# in the source workbook - replicate by copying that cell format over.
$ws.Range("E5").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null

$ws.Range("E6").Copy() | Out-Null
$ws.Range("H6").PasteSpecial(-4122) | Out-Null

$ws.Range("E6").Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("I4").Select() | Out-Null
